$wb = $excel.ActiveWorkbook

# 1. Rename the second sheet "Referencias" -> "Referentes"
$wsObra = $wb.Worksheets.Item("Obra")
$wsRef = $wb.Worksheets.Item("Referencias")
$wsRef.Name = "Referentes"

# 2. Sheet "Obra": the "Laminas de paisajes latinoamericanos" row (Id 3, row 4) is
#    removed, and the row below it ("Zocalo de la tragedia", Id 4, row 5) shifts up
#    to become row 4, renumbered to Id 3.
$wsObra.Rows.Item(4).Delete()
$wsObra.Range("A4").Value = 3

# 3. Sheet "Referentes": row 4 (Id 3) gains the "Laminas de paisajes
#    latinoamericanos" / "laminas-paisajes.jpg" data that was dropped from "Obra".
$wsRef.Range("B3").Copy()
$wsRef.Range("B4").PasteSpecial(-4122)
$wsRef.Range("A2").Copy()
$wsRef.Range("E4").PasteSpecial(-4122)
$wsRef.Range("B4").Value = "Láminas de paisajes latinoamericanos"
$wsRef.Range("E4").Value = "laminas-paisajes.jpg"
